$d = $word.ActiveDocument

# Remove all comments from the document, keeping the underlying text intact.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}
$d.Save()
